# JP tenth anni — backfill TW names for rows 155-162, CN names for rows 166-170,
# and append two new servants (rows 179-180).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Backfill missing name_TW (column F) values
$ws.Range("F155").Value = "洛庫斯塔"
$ws.Range("F156").Value = "瑟坦特"
$ws.Range("F157").Value = "難敵"
$ws.Range("F158").Value = "鈴鹿御前〔暑假〕"
$ws.Range("F159").Value = "克洛伊．馮．愛因茲貝倫"
$ws.Range("F160").Value = "諾克娜蕾雅．雅蘭杜"
$ws.Range("F161").Value = "ＵＤＫ－巴格斯特"
$ws.Range("F162").Value = "凱特．庫．米可科爾"

# Backfill missing name_CN (column E) values
$ws.Range("E166").Value = "由井正雪"
$ws.Range("E167").Value = "宮本伊織"
$ws.Range("E168").Value = "耀星哈桑"
$ws.Range("E169").Value = "亞歷山德羅·迪·卡利奧斯特羅"
$ws.Range("E170").Value = "靜希草十郎"

# Append new servants
$ws.Range("A179").Value = 434
$ws.Range("B179").Value = 4
$ws.Range("C179").Value = "Saber"
$ws.Range("D179").Value = "黒姫"
# name_CN / name_TW are not published yet for this servant, but the column
# is still a (blank) text cell for this row, same as every other row —
# write an explicit empty-text value via the leading-apostrophe text marker,
# then strip the resulting quote-prefix formatting so no stray style sticks.
$ws.Range("E179").Value = "'"
$ws.Range("E179").Style = "Normal"
$ws.Range("F179").Value = "'"
$ws.Range("F179").Style = "Normal"

$ws.Range("A180").Value = 439
$ws.Range("B180").Value = 4
$ws.Range("C180").Value = "Lancer"
$ws.Range("D180").Value = "アショカ王"
$ws.Range("E180").Value = "'"
$ws.Range("E180").Style = "Normal"
$ws.Range("F180").Value = "'"
$ws.Range("F180").Style = "Normal"
